# Update Kenya household sanitation facilities data:
#  - add 2016 survey values alongside the existing 2013 values on the
#    "Data" sheet (one new row per district)
#  - add the corresponding "2016" column to the "Data-wide-value" sheet
#  - bump the "Source:" note on the "Notes" sheet

$wb = $excel.ActiveWorkbook

$districts = @(
    @{ id = 'd18974'; name = 'Baringo'; v2013 = 41.0488245931284; v2016 = 71.3 },
    @{ id = 'd18975'; name = 'Bomet'; v2013 = 69.0179806362379; v2016 = 30.6 },
    @{ id = 'd18987'; name = 'Bungoma'; v2013 = 71.5288035450517; v2016 = 39.3 },
    @{ id = 'd18988'; name = 'Busia'; v2013 = 60.6355382619974; v2016 = 61.6 },
    @{ id = 'd18976'; name = 'Elgeyo-Marakwet'; v2013 = 48.4536082474227; v2016 = 72.9 },
    @{ id = 'd18955'; name = 'Embu'; v2013 = 65.5277145026576; v2016 = 98.4 },
    @{ id = 'd18965'; name = 'Garissa'; v2013 = 31.0344827586207; v2016 = 45.1 },
    @{ id = 'd18968'; name = 'Homa Bay'; v2013 = 40.6886517943744; v2016 = 25.1 },
    @{ id = 'd18956'; name = 'Isiolo'; v2013 = 42.1725239616613; v2016 = 70.7 },
    @{ id = 'd18991'; name = 'Kajiado'; v2013 = 60.9798270893372; v2016 = 84.3 },
    @{ id = 'd18989'; name = 'Kakamega'; v2013 = 83.4973292100084; v2016 = 32.4 },
    @{ id = 'd18977'; name = 'Kericho'; v2013 = 64.3445218492717; v2016 = 79 },
    @{ id = 'd18943'; name = 'Kiambu'; v2013 = 80.4560954816709; v2016 = 90.3 },
    @{ id = 'd18949'; name = 'Kilifi'; v2013 = 47.9479479479479; v2016 = 70 },
    @{ id = 'd18944'; name = 'Kirinyaga'; v2013 = 82.8793774319066; v2016 = 87 },
    @{ id = 'd18969'; name = 'Kisii'; v2013 = 63.9591836734694; v2016 = 41.3 },
    @{ id = 'd18970'; name = 'Kisumu'; v2013 = 56.2169312169312; v2016 = 96.3 },
    @{ id = 'd18957'; name = 'Kitui'; v2013 = 54.014598540146; v2016 = 56.8 },
    @{ id = 'd18950'; name = 'Kwale'; v2013 = 34.0163934426229; v2016 = 33.7 },
    @{ id = 'd18978'; name = 'Laikipia'; v2013 = 72.2599418040737; v2016 = 42.3 },
    @{ id = 'd18951'; name = 'Lamu'; v2013 = 57.2072072072072; v2016 = 70 },
    @{ id = 'd18958'; name = 'Machakos'; v2013 = 62.7599243856333; v2016 = 83.5 },
    @{ id = 'd18959'; name = 'Makueni'; v2013 = 56.3538873994638; v2016 = 88 },
    @{ id = 'd18966'; name = 'Mandera'; v2013 = 16.3346613545817; v2016 = 40.9 },
    @{ id = 'd18960'; name = 'Marsabit'; v2013 = 25.8347978910369; v2016 = 26 },
    @{ id = 'd18961'; name = 'Meru'; v2013 = 78.0976220275344; v2016 = 48.7 },
    @{ id = 'd18971'; name = 'Migori'; v2013 = 51.3588463671658; v2016 = 37.9 },
    @{ id = 'd18952'; name = 'Mombasa'; v2013 = 81.2058057312988; v2016 = 86.3 },
    @{ id = 'd18946'; name = 'Murang''a'; v2013 = 69.1044192412984; v2016 = 65.6 },
    @{ id = 'd18964'; name = 'Nairobi'; v2013 = 86.8426395939086; v2016 = 91.9 },
    @{ id = 'd18979'; name = 'Nakuru'; v2013 = 76.8667642752562; v2016 = 62.1 },
    @{ id = 'd18980'; name = 'Nandi'; v2013 = 72.680077871512; v2016 = 71 },
    @{ id = 'd18981'; name = 'Narok'; v2013 = 36.7021276595745; v2016 = 42.1 },
    @{ id = 'd18962'; name = 'Nithi'; v2013 = 60.5855855855856; v2016 = 96.3 },
    @{ id = 'd18972'; name = 'Nyamira'; v2013 = 66.8702290076336; v2016 = 40.4 },
    @{ id = 'd18947'; name = 'Nyandarua'; v2013 = 74.0792216817234; v2016 = 81.8 },
    @{ id = 'd18948'; name = 'Nyeri'; v2013 = 74.3182944967774; v2016 = 55.2 },
    @{ id = 'd18982'; name = 'Samburu'; v2013 = 21.5644820295983; v2016 = 31.5 },
    @{ id = 'd18973'; name = 'Siaya'; v2013 = 47.8894472361809; v2016 = 43.7 },
    @{ id = 'd18953'; name = 'Taita Taveta'; v2013 = 67.7918424753868; v2016 = 98.8 },
    @{ id = 'd18954'; name = 'Tana River'; v2013 = 22.9957805907173; v2016 = 60.6 },
    @{ id = 'd18983'; name = 'Trans-Nzoia'; v2013 = 74.9706227967097; v2016 = 36.1 },
    @{ id = 'd18984'; name = 'Turkana'; v2013 = 12.1753246753247; v2016 = 31.9 },
    @{ id = 'd18985'; name = 'Uasin Gishu'; v2013 = 78.596144340089; v2016 = 83.2 },
    @{ id = 'd18990'; name = 'Vihiga'; v2013 = 85.6447688564477; v2016 = 43.2 },
    @{ id = 'd18967'; name = 'Wajir'; v2013 = 7.44920993227991; v2016 = 6.2 },
    @{ id = 'd18986'; name = 'West Pokot'; v2013 = 26.7590618336887; v2016 = 27.3 }
)

# --- "Notes" sheet: update source reference -------------------------------
$notes = $wb.Worksheets.Item("Notes")
$notes.Range("A4").Value2 = "Source: 38"

# --- "Data" sheet: insert a 2016 row right after each district's 2013 row -
$data = $wb.Worksheets.Item("Data")

$row = 2
foreach ($d in $districts) {
    # existing 2013 row (row $row) is left untouched

    # new 2016 row, inserted directly below it
    $row = $row + 1
    $data.Rows.Item($row).Insert()
    $data.Cells.Item($row, 1).Value2 = $d.id
    $data.Cells.Item($row, 2).Value2 = $d.name
    $data.Cells.Item($row, 3).Value2 = 2016
    $data.Cells.Item($row, 4).Value2 = $d.v2016
    $row = $row + 1
}

# --- "Data-wide-value" sheet: add a "2016" column --------------------------
$wide = $wb.Worksheets.Item("Data-wide-value")
$wide.Range("C1").NumberFormat = "@"
$wide.Range("C1").Value2 = "2016"

$row = 2
foreach ($d in $districts) {
    $wide.Cells.Item($row, 3).Value2 = $d.v2016
    $row = $row + 1
}
